$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(3, 8).Value = 2.06
$ws.Cells.Item(4, 6).Value = 1.61
$ws.Cells.Item(4, 7).Value = 1.67
$ws.Cells.Item(4, 20).Value = 1.72
$ws.Cells.Item(4, 21).Value = 2.18
$ws.Cells.Item(5, 6).Value = 1.42
$ws.Cells.Item(5, 7).Value = 1.43
$ws.Cells.Item(5, 8).Value = 8.4
$ws.Cells.Item(5, 9).Value = 9.199999999999999
$ws.Cells.Item(5, 10).Value = 5.3
$ws.Cells.Item(5, 11).Value = 5.6
$ws.Cells.Item(5, 17).Value = 1.8
$ws.Cells.Item(5, 20).Value = 2.12
$ws.Cells.Item(5, 24).Value = 23
$ws.Cells.Item(5, 31).Value = 150
$ws.Cells.Item(5, 34).Value = 55
$ws.Cells.Item(5, 38).Value = 980
$ws.Cells.Item(5, 39).Value = 180
$ws.Cells.Item(5, 41).Value = 220
$ws.Cells.Item(6, 6).Value = 5.1
$ws.Cells.Item(6, 8).Value = 1.79
$ws.Cells.Item(6, 21).Value = 1.93
$ws.Cells.Item(7, 11).Value = 8.6
$ws.Cells.Item(7, 17).Value = 1.42
$ws.Cells.Item(7, 19).Value = 2.04
$ws.Cells.Item(7, 20).Value = 2.02
$ws.Cells.Item(7, 21).Value = 1.9
$ws.Cells.Item(8, 8).Value = 4.7
$ws.Cells.Item(8, 14).Value = 5.2
$ws.Cells.Item(8, 17).Value = 1.62
$ws.Cells.Item(8, 19).Value = 2.54
$ws.Cells.Item(8, 21).Value = 2.48
$ws.Cells.Item(9, 7).Value = 1.6
$ws.Cells.Item(9, 8).Value = 6.4
$ws.Cells.Item(9, 9).Value = 7
$ws.Cells.Item(9, 20).Value = 1.83
$ws.Cells.Item(9, 21).Value = 2.1
$ws.Cells.Item(10, 9).Value = 7.6
$ws.Cells.Item(10, 16).Value = 2.34
$ws.Cells.Item(10, 17).Value = 1.68
$ws.Cells.Item(11, 9).Value = 2.46
$ws.Cells.Item(12, 6).Value = 1.56
$ws.Cells.Item(12, 7).Value = 1.59
$ws.Cells.Item(12, 8).Value = 6
$ws.Cells.Item(12, 9).Value = 6.6
$ws.Cells.Item(12, 16).Value = 2.32
$ws.Cells.Item(12, 17).Value = 1.7
$ws.Cells.Item(12, 20).Value = 1.81
$ws.Cells.Item(12, 21).Value = 2.16
$ws.Cells.Item(12, 24).Value = 29
$ws.Cells.Item(12, 25).Value = 34
$ws.Cells.Item(12, 26).Value = 190
$ws.Cells.Item(12, 29).Value = 11
$ws.Cells.Item(12, 30).Value = 34
$ws.Cells.Item(12, 34).Value = 29
$ws.Cells.Item(12, 40).Value = 7.2
$ws.Cells.Item(13, 6).Value = 5
$ws.Cells.Item(13, 7).Value = 5.9
$ws.Cells.Item(13, 8).Value = 1.76
$ws.Cells.Item(13, 9).Value = 1.79
$ws.Cells.Item(13, 10).Value = 3.95
$ws.Cells.Item(13, 16).Value = 2.14
$ws.Cells.Item(13, 17).Value = 1.71
$ws.Cells.Item(14, 6).Value = 2.26
$ws.Cells.Item(14, 7).Value = 2.5
$ws.Cells.Item(14, 8).Value = 3.35
$ws.Cells.Item(14, 9).Value = 3.65
$ws.Cells.Item(14, 10).Value = 3.35
$ws.Cells.Item(14, 16).Value = 1.79
$ws.Cells.Item(14, 17).Value = 1.9
$ws.Cells.Item(15, 9).Value = 16.5
$ws.Cells.Item(15, 16).Value = 2.8
$ws.Cells.Item(15, 17).Value = 1.45
$ws.Cells.Item(16, 6).Value = 1.64
$ws.Cells.Item(16, 10).Value = 3.9
$ws.Cells.Item(18, 6).Value = 1.89
$ws.Cells.Item(18, 7).Value = 1.97
$ws.Cells.Item(18, 10).Value = 4.1
$ws.Cells.Item(18, 16).Value = 2.82
$ws.Cells.Item(18, 17).Value = 1.44
$ws.Cells.Item(19, 6).Value = 2.58
$ws.Cells.Item(21, 16).Value = 1.84
$ws.Cells.Item(21, 17).Value = 1.98
$ws.Cells.Item(22, 11).Value = 5.1
$ws.Cells.Item(22, 13).Value = 1.03
$ws.Cells.Item(22, 21).Value = 2.14
$ws.Cells.Item(22, 33).Value = 1000
$ws.Cells.Item(24, 20).Value = 1.6
$ws.Cells.Item(24, 21).Value = 2.52
$ws.Cells.Item(25, 6).Value = 1.77
$ws.Cells.Item(25, 16).Value = 1.97
$ws.Cells.Item(25, 17).Value = 1.92
$ws.Cells.Item(25, 20).Value = 1.9
$ws.Cells.Item(25, 29).Value = 8.6
$ws.Cells.Item(25, 33).Value = 12
$ws.Cells.Item(26, 6).Value = 4.9
$ws.Cells.Item(26, 8).Value = 1.76
$ws.Cells.Item(26, 26).Value = 11.5
$ws.Cells.Item(27, 17).Value = 1.67
$ws.Cells.Item(28, 17).Value = 1.74
$ws.Cells.Item(28, 27).Value = 290
$ws.Cells.Item(28, 38).Value = 1000
$ws.Cells.Item(29, 15).Value = 1.22
$ws.Cells.Item(29, 17).Value = 1.66
$ws.Cells.Item(29, 19).Value = 2.64
$ws.Cells.Item(29, 24).Value = 28
$ws.Cells.Item(30, 6).Value = 2.26
$ws.Cells.Item(30, 7).Value = 2.38
$ws.Cells.Item(30, 8).Value = 3.3
$ws.Cells.Item(30, 9).Value = 3.65
$ws.Cells.Item(30, 10).Value = 3.25
$ws.Cells.Item(30, 11).Value = 3.65
$ws.Cells.Item(30, 16).Value = 1.72
$ws.Cells.Item(30, 17).Value = 2.1
$ws.Cells.Item(31, 9).Value = 3.65
$ws.Cells.Item(31, 10).Value = 3.45
$ws.Cells.Item(31, 16).Value = 2
$ws.Cells.Item(31, 17).Value = 1.81
$ws.Cells.Item(33, 6).Value = 1.45
$ws.Cells.Item(33, 7).Value = 1.64
$ws.Cells.Item(33, 8).Value = 6.4
$ws.Cells.Item(33, 9).Value = 1000
$ws.Cells.Item(33, 10).Value = 3.7
$ws.Cells.Item(33, 11).Value = 5.2
$ws.Cells.Item(33, 16).Value = 2.3
$ws.Cells.Item(33, 17).Value = 1.01
$ws.Cells.Item(34, 6).Value = 4.2
$ws.Cells.Item(34, 8).Value = 1.89
$ws.Cells.Item(34, 9).Value = 2.02
$ws.Cells.Item(34, 10).Value = 3.65
$ws.Cells.Item(34, 11).Value = 3.9
$ws.Cells.Item(34, 16).Value = 1.94
$ws.Cells.Item(34, 17).Value = 1.89
$ws.Cells.Item(35, 6).Value = 2.46
$ws.Cells.Item(35, 10).Value = 3.45
$ws.Cells.Item(35, 11).Value = 3.6
$ws.Cells.Item(35, 16).Value = 1.95
$ws.Cells.Item(35, 17).Value = 1.84
$ws.Cells.Item(36, 7).Value = 8.6
$ws.Cells.Item(36, 9).Value = 1.52
$ws.Cells.Item(36, 11).Value = 5.6
$ws.Cells.Item(36, 16).Value = 2.26
$ws.Cells.Item(36, 17).Value = 1.64
$ws.Cells.Item(37, 6).Value = 1.71
$ws.Cells.Item(37, 7).Value = 1.81
$ws.Cells.Item(37, 9).Value = 6.4
$ws.Cells.Item(37, 16).Value = 1.79
$ws.Cells.Item(38, 8).Value = 2.3
$ws.Cells.Item(38, 15).Value = 1.33
$ws.Cells.Item(38, 21).Value = 2.16
$ws.Cells.Item(39, 9).Value = 3.1
